# Fix salvage value calculation issue.
# Salvage Value should be 20% of Initial Investment; several rows had a
# stale, hard-coded Salvage Value (and, for two machines, a stale
# Economic Life). This corrects those inputs and the values that were
# derived from them (Use Cost on "Summary Costs", Maintenance & Repairs
# on "Operating Costs", and the depreciation/interest/insurance/taxes
# figures on "Fixed Costs").

$wb = $excel.ActiveWorkbook

$wsSummary   = $wb.Worksheets.Item("Summary Costs")
$wsOperating = $wb.Worksheets.Item("Operating Costs")
$wsFixed     = $wb.Worksheets.Item("Fixed Costs")

# ---------------------------------------------------------------------
# "Summary Costs" sheet: Salvage Value (D), Economic Life (E), Use Cost (I)
# ---------------------------------------------------------------------
$wsSummary.Range("D2").Value = 40000
$wsSummary.Range("E2").Value = 10
$wsSummary.Range("I2").Value = 56.5070150452

$wsSummary.Range("D3").Value = 74800
$wsSummary.Range("E3").Value = 10
$wsSummary.Range("I3").Value = 87.3352722135

$wsSummary.Range("D4").Value = 51000
$wsSummary.Range("I4").Value = 74.7452978886

$wsSummary.Range("D5").Value = 23611.2
$wsSummary.Range("I5").Value = 43.2408100045

$wsSummary.Range("D7").Value = 27589.8
$wsSummary.Range("I7").Value = 42.7877603491

$wsSummary.Range("D8").Value = 54000
$wsSummary.Range("I8").Value = 56.767622612

$wsSummary.Range("D9").Value = 101000
$wsSummary.Range("I9").Value = 99.513712625

# ---------------------------------------------------------------------
# "Operating Costs" sheet: Maintenance and Repairs ($/hr) (C)
# ---------------------------------------------------------------------
$wsOperating.Range("C2").Value = 14.4813059505
$wsOperating.Range("C3").Value = 27.0800421274
$wsOperating.Range("C4").Value = 24.6182201159
$wsOperating.Range("C5").Value = 10.4055769231
$wsOperating.Range("C7").Value = 12.485510466
$wsOperating.Range("C8").Value = 7.61538461538
$wsOperating.Range("C9").Value = 25.3219373219

# ---------------------------------------------------------------------
# "Fixed Costs" sheet: Annual Depreciation (C), Annual Interest (D),
# Annual Insurance (E), Annual Taxes (F)
# ---------------------------------------------------------------------
$wsFixed.Range("C2").Value = 16000
$wsFixed.Range("D2").Value = 15360
$wsFixed.Range("E2").Value = 3840
$wsFixed.Range("F2").Value = 3840

$wsFixed.Range("C3").Value = 29920
$wsFixed.Range("D3").Value = 28723.2
$wsFixed.Range("E3").Value = 7180.8
$wsFixed.Range("F3").Value = 7180.8

$wsFixed.Range("C4").Value = 27200
$wsFixed.Range("D4").Value = 19992
$wsFixed.Range("E4").Value = 4998
$wsFixed.Range("F4").Value = 4998

$wsFixed.Range("C5").Value = 11805.6
$wsFixed.Range("D5").Value = 9208.368
$wsFixed.Range("E5").Value = 2302.092
$wsFixed.Range("F5").Value = 2302.092

$wsFixed.Range("C7").Value = 13794.9
$wsFixed.Range("D7").Value = 10760.022
$wsFixed.Range("E7").Value = 2690.0055
$wsFixed.Range("F7").Value = 2690.0055

$wsFixed.Range("C8").Value = 10800
$wsFixed.Range("D8").Value = 20088
$wsFixed.Range("E8").Value = 5022
$wsFixed.Range("F8").Value = 5022

$wsFixed.Range("C9").Value = 40400
$wsFixed.Range("D9").Value = 38784
$wsFixed.Range("E9").Value = 9696
$wsFixed.Range("F9").Value = 9696
